$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old inline string content that is being replaced
$ws.Range("A1").Value = $null
$ws.Range("B1").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("F10").Value = $null

# Set the new numeric diagonal values
$ws.Range("A1").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D4").Value = 4
$ws.Range("E5").Value = 5
$ws.Range("F6").Value = 6
$ws.Range("G7").Value = 7
$ws.Range("H8").Value = 8
$ws.Range("I9").Value = 9
$ws.Range("J10").Value = 10
